# Update generated data values as reflected in the "gh-pages" output regeneration.
$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 903 -> 907, F5 536 -> 537
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 907
$wsExhibit.Range("F5").Value = 537

# Sheet "全部类型": F4 903 -> 907, F6 536 -> 537
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 907
$wsAll.Range("F6").Value = 537
